$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/percentage updates (safe to assign directly; Excel will not
# reinterpret these as numbers because of letters, "%", spaces, or multiple
# "." separators).
$ws.Range("D2").Value = '34.503.39'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.811.41'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("E6").Value = '  +2.98%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +5.91%  '
$ws.Range("E9").Value = '  -4.20%  '
$ws.Range("E10").Value = '  -2.50%  '
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("E13").Value = '  -3.81%  '
$ws.Range("D14").Value = '1.806.20'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("D16").Value = '34.467.97'
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").Value = '0.0₃0774'
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("E21").Value = '  -2.52%  '
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("E24").Value = '  +3.91%  '
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("E27").Value = '  +4.86%  '
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("E32").Value = '  -2.85%  '
$ws.Range("E33").Value = '  -4.45%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").Value = '1.357.89'
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("E36").Value = '  -4.44%  '
$ws.Range("E37").Value = '  -0.64%  '
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("E39").Value = '  -4.93%  '
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("E45").Value = '  +1.73%  '
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("D47").Value = '1.973.49'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("E50").Value = '  -2.14%  '
$ws.Range("D51").Value = '0.0₆0124'
$ws.Range("E51").Value = '  -3.64%  '

# Purely-numeric-looking price text (e.g. "225.63") must be forced to stay
# text, matching the source data (column D is always text in this sheet).
# Temporarily mark the cell as Text, assign the string, then restore the
# default "Normal" style so no stray number-format sticks to the cell.
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '225.63'
$ws.Range("D5").Style = 'Normal'

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '38.14'
$ws.Range("D8").Style = 'Normal'

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.0676'
$ws.Range("D10").Style = 'Normal'

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '11.26'
$ws.Range("D13").Style = 'Normal'

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.633'
$ws.Range("D15").Style = 'Normal'

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '68.39'
$ws.Range("D18").Style = 'Normal'

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '243.33'
$ws.Range("D19").Style = 'Normal'

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '11.22'
$ws.Range("D21").Style = 'Normal'

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '170.68'
$ws.Range("D25").Style = 'Normal'

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '7.78'
$ws.Range("D26").Style = 'Normal'

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '17.72'
$ws.Range("D27").Style = 'Normal'

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '1.83'
$ws.Range("D34").Style = 'Normal'

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '2.45'
$ws.Range("D40").Style = 'Normal'

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.956'
$ws.Range("D41").Style = 'Normal'

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '1.21'
$ws.Range("D42").Style = 'Normal'

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '81.92'
$ws.Range("D43").Style = 'Normal'

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '102.61'
$ws.Range("D50").Style = 'Normal'
